$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DUI traffic data values in column E (year values)
$ws.Range("E15").Value = 2020
$ws.Range("E59").Value = 2019
$ws.Range("E60").Value = 2019
$ws.Range("E61").Value = 2019
$ws.Range("E62").Value = 2019

# Update sheet view: scroll the window so row 43 is at the top, then select G16
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G16").Select()
